# Adding excel files and creating get available usd
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "safety_orders": insert a new safety order row (new row 2)
# and refresh the recalculated values for all rows.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("safety_orders")

$ws1.Rows.Item(2).Insert()
$ws1.Rows.Item(2).ClearFormats()

$ws1.Cells.Item(2,1).Value = 2
$ws1.Cells.Item(2,2).Value = 3.328
$ws1.Cells.Item(2,3).Value = 0.125
$ws1.Cells.Item(2,4).Value = 61.90198176
$ws1.Cells.Item(2,5).Value = 62.55127638
$ws1.Cells.Item(2,6).Value = 63.17678914380001
$ws1.Cells.Item(2,7).Value = 2.017841364014805

$ws1.Cells.Item(3,1).Value = 3
$ws1.Cells.Item(3,2).Value = 6.49168
$ws1.Cells.Item(3,3).Value = 0.3125
$ws1.Cells.Item(3,4).Value = 59.8761825456
$ws1.Cells.Item(3,5).Value = 61.2137294628
$ws1.Cells.Item(3,6).Value = 61.825866757428
$ws1.Cells.Item(3,7).Value = 3.153508901828306

$ws1.Cells.Item(4,1).Value = 4
$ws1.Cells.Item(4,2).Value = 11.4270208
$ws1.Cells.Item(4,3).Value = 0.78125
$ws1.Cells.Item(4,4).Value = 56.715935771136
$ws1.Cells.Item(4,5).Value = 58.964832616968
$ws1.Cells.Item(4,6).Value = 59.55448094313768
$ws1.Cells.Item(4,7).Value = 4.766299910684991

$ws1.Cells.Item(5,1).Value = 5
$ws1.Cells.Item(5,2).Value = 19.12615245
$ws1.Cells.Item(5,3).Value = 1.953125
$ws1.Cells.Item(5,4).Value = 51.7859508016915
$ws1.Cells.Item(5,5).Value = 55.37539170932975
$ws1.Cells.Item(5,6).Value = 55.92914562642304
$ws1.Cells.Item(5,7).Value = 7.407935126357701

$ws1.Cells.Item(6,1).Value = 6
$ws1.Cells.Item(6,2).Value = 31.13679782
$ws1.Cells.Item(6,3).Value = 4.8828125
$ws1.Cells.Item(6,4).Value = 44.0951742519194
$ws1.Cells.Item(6,5).Value = 49.73528298062458
$ws1.Cells.Item(6,6).Value = 50.23263581043082
$ws1.Cells.Item(6,7).Value = 12.21807587735019

$ws1.Cells.Item(7,1).Value = 7
$ws1.Cells.Item(7,2).Value = 49.8734046
$ws1.Cells.Item(7,3).Value = 12.20703125
$ws1.Cells.Item(7,4).Value = 32.097562832482
$ws1.Cells.Item(7,5).Value = 40.91642290655329
$ws1.Cells.Item(7,6).Value = 41.32558713561882
$ws1.Cells.Item(7,7).Value = 22.33005007975585

# ---------------------------------------------------------------
# Sheet "open_buy_orders": replace the single remaining open buy
# order and drop the now-filled second order row.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("open_buy_orders")

$ws2.Cells.Item(2,1).Value = "OWPVJY-PJYD4-FKO65H"
$ws2.Cells.Item(2,2).Value = 63.832

$ws2.Rows.Item(3).Delete()

# ---------------------------------------------------------------
# Sheet "open_sell_orders": replace the open sell order txid.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("open_sell_orders")

$ws3.Cells.Item(2,1).Value = "OHLGAU-Q63PC-F7C3CX"
